$d = $word.ActiveDocument

# The paragraph that ends "...enforce higher code quality." is the 4th
# paragraph in the document; the new "Features " heading paragraph must be
# inserted right after it and before the trailing empty paragraph.
$target = $d.Paragraphs(4)
$target.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs(5)

# Build the new paragraph (bold, justified, size 24 half-points = 12pt,
# text "Features ") as a raw OOXML package fragment so the run formatting
# (w:b / w:bCs) is applied exactly, including on the paragraph mark.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Features </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xml) | Out-Null
